$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03835166666666667
$ws.Range("H2").Value = 0.115055
$ws.Range("I2").Value = 0.0181239951898282
$ws.Range("J2").Value = 0.0181239951898282
$ws.Range("M2").Value = 19.827687
$ws.Range("N2").Value = 59.483061
$ws.Range("O2").Value = 0.1538389073329896
$ws.Range("P2").Value = 0.1538389073329896
$ws.Range("Q2").Value = 0.760424842595
$ws.Range("R2").Value = 6.843823583355
$ws.Range("S2").Value = 0.00278817561651153
$ws.Range("T2").Value = 0.002788175616511529

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03835166666666667
$ws.Range("H3").Value = 0.115055
$ws.Range("I3").Value = 0.0181239951898282
$ws.Range("J3").Value = 0.0181239951898282
$ws.Range("O3").Value = 0.6604253914664442
$ws.Range("P3").Value = 0.6604253914664441
$ws.Range("Q3").Value = 3.264478947868334
$ws.Range("R3").Value = 29.380310530815
$ws.Range("S3").Value = 0.01196954661817824
$ws.Range("T3").Value = 0.01196954661817824

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.03835166666666667
$ws.Range("H4").Value = 0.115055
$ws.Range("I4").Value = 0.0181239951898282
$ws.Range("J4").Value = 0.0181239951898282
$ws.Range("M4").Value = 23.93873833333333
$ws.Range("N4").Value = 71.816215
$ws.Range("O4").Value = 0.1857357012005663
$ws.Range("P4").Value = 0.1857357012005663
$ws.Range("Q4").Value = 0.9180905129805556
$ws.Range("R4").Value = 8.262814616825001
$ws.Range("S4").Value = 0.003366272955138431
$ws.Range("T4").Value = 0.00336627295513843

$ws.Range("I5").Value = 0.3727881574250648
$ws.Range("J5").Value = 0.3727881574250648
$ws.Range("M5").Value = 19.827687
$ws.Range("N5").Value = 59.483061
$ws.Range("O5").Value = 0.1538389073329896
$ws.Range("P5").Value = 0.1538389073329896
$ws.Range("Q5").Value = 15.640998188431
$ws.Range("R5").Value = 140.768983695879
$ws.Range("S5").Value = 0.05734932280495049
$ws.Range("T5").Value = 0.05734932280495048

$ws.Range("I6").Value = 0.3727881574250648
$ws.Range("J6").Value = 0.3727881574250648
$ws.Range("O6").Value = 0.6604253914664442
$ws.Range("P6").Value = 0.6604253914664441
$ws.Range("S6").Value = 0.2461987648015028
$ws.Range("T6").Value = 0.2461987648015028

$ws.Range("I7").Value = 0.3727881574250648
$ws.Range("J7").Value = 0.3727881574250648
$ws.Range("M7").Value = 23.93873833333333
$ws.Range("N7").Value = 71.816215
$ws.Range("O7").Value = 0.1857357012005663
$ws.Range("P7").Value = 0.1857357012005663
$ws.Range("Q7").Value = 18.88398595887611
$ws.Range("R7").Value = 169.955873629885
$ws.Range("S7").Value = 0.06924006981861151
$ws.Range("T7").Value = 0.0692400698186115

$ws.Range("G8").Value = 1.288873333333333
$ws.Range("H8").Value = 3.86662
$ws.Range("I8").Value = 0.6090878473851071
$ws.Range("J8").Value = 0.609087847385107
$ws.Range("M8").Value = 19.827687
$ws.Range("N8").Value = 59.483061
$ws.Range("O8").Value = 0.1538389073329896
$ws.Range("P8").Value = 0.1538389073329896
$ws.Range("Q8").Value = 25.55537703598
$ws.Range("R8").Value = 229.99839332382
$ws.Range("S8").Value = 0.09370140891152763
$ws.Range("T8").Value = 0.09370140891152759

$ws.Range("G9").Value = 1.288873333333333
$ws.Range("H9").Value = 3.86662
$ws.Range("I9").Value = 0.6090878473851071
$ws.Range("J9").Value = 0.609087847385107
$ws.Range("O9").Value = 0.6604253914664442
$ws.Range("P9").Value = 0.6604253914664441
$ws.Range("Q9").Value = 109.7083967616067
$ws.Range("R9").Value = 987.37557085446
$ws.Range("S9").Value = 0.4022570800467632
$ws.Range("T9").Value = 0.402257080046763

$ws.Range("G10").Value = 1.288873333333333
$ws.Range("H10").Value = 3.86662
$ws.Range("I10").Value = 0.6090878473851071
$ws.Range("J10").Value = 0.609087847385107
$ws.Range("M10").Value = 23.93873833333333
$ws.Range("N10").Value = 71.816215
$ws.Range("O10").Value = 0.1857357012005663
$ws.Range("P10").Value = 0.1857357012005663
$ws.Range("Q10").Value = 30.85400147147778
$ws.Range("R10").Value = 277.6860132433
$ws.Range("S10").Value = 0.1131293584268164
$ws.Range("T10").Value = 0.1131293584268164
